$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 updates
$ws.Range("A16").Value = 112178652
$ws.Range("B16").Value = 90812
$ws.Range("D16").Value = "LC"
$ws.Range("E16").Value = 4366
$ws.Range("F16").Value = "Skarp dropptaggsvamp"
$ws.Range("G16").Value = "Hydnellum peckii"
$ws.Range("H16").Value = "Banker"
$ws.Range("Q16").Value = 618476
$ws.Range("R16").Value = 6905002

# Row 17 updates
$ws.Range("B17").Value = 86357

# Row 18 updates
$ws.Range("A18").Value = 112178654
$ws.Range("B18").Value = 89820
$ws.Range("D18").Value = "NT"
$ws.Range("E18").Value = 658
$ws.Range("F18").Value = "Rosenticka"
$ws.Range("G18").Value = "Rhodofomes roseus"
$ws.Range("H18").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("Q18").Value = 618387
$ws.Range("R18").Value = 6904851
